$wb = $excel.ActiveWorkbook

# "Active" sheet: the task "scroll horizontal and vertical when zoomed in" (Id 4)
# is complete, so remove its row (row 3); subsequent rows shift up.
$wsActive = $wb.Worksheets.Item("Active")
$wsActive.Rows.Item(3).Delete()

# "Inactive" sheet: insert that same task at the top (row 2), marked Done,
# with a Done date of 8/10/2018. Existing rows shift down by one.
$wsInactive = $wb.Worksheets.Item("Inactive")
$wsInactive.Rows.Item(2).Insert()
$wsInactive.Range("A2").Value = 4
$wsInactive.Range("B2").Value = "scroll horizontal and vertical when zoomed in"
$wsInactive.Range("C2").Value = "Done"
$wsInactive.Range("D2").Value = "Task"
$wsInactive.Range("E2").Value = "'8/9/2018"
$wsInactive.Range("F2").Value = "'8/10/2018"
$wsInactive.Range("A2:F2").Style = "Normal"
